$wb = $excel.ActiveWorkbook

# --- Update test data values (Pro Lite panel ticket references) ---

# Portugal: ticket reference gains an additional linked task
$wsPortugal = $wb.Worksheets.Item("Portugal")
$wsPortugal.Range("B4").Value = "NGC-3479/T2407/T2508"

# Belgium: ticket reference renumbered
$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsBelgium.Range("B4").Value = "NGC-3478/T2263"

# Czech: ticket reference renumbered
$wsCzech = $wb.Worksheets.Item("Czech")
$wsCzech.Range("B4").Value = "NGC-3477/T1732"

# --- Update sheet selections / active view state ---

$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Range("A7").Select()

$wsCzech.Range("B3").Select()

$wsPortugal.Range("A8").Select()

# Belgium becomes the active tab, selection moves to A9
$wsBelgium.Activate()
$wsBelgium.Range("A9").Select()
